$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / metadata changes (wear -> wide rebrand) ---
$ws.Range("B1").Value = "http://purl.org/wide/"
$ws.Range("B3").Value = "wide"
$ws.Range("C3").Value = "http://purl.org/wide/"
$ws.Range("B8").Value = "WInd energy moDEls (WIDE) Taxonomy"
$ws.Range("B9").Value = "A classification of models used in wind energy"

# --- Rows 17-54: replace taxonomy term identifiers/labels/broader refs ---
$ws.Range("A17").Value = "wide:MeteorologicalModels"
$ws.Range("B17").Value = "Meteorological Models"
$ws.Range("D17").Value = ""
$ws.Range("F17").Value = ""

$ws.Range("A18").Value = "wide:OceanographicModels"
$ws.Range("B18").Value = "Oceanographic Models"
$ws.Range("D18").Value = ""
$ws.Range("F18").Value = ""

$ws.Range("A19").Value = "wide:HydrodynamicModels"
$ws.Range("B19").Value = "Hydrodynamic Models"
$ws.Range("D19").Value = ""
$ws.Range("F19").Value = "wide:OceanographicModels"

$ws.Range("A20").Value = "wide:WaveModels"
$ws.Range("B20").Value = "Wave Models"
$ws.Range("D20").Value = ""
$ws.Range("F20").Value = "wide:OceanographicModels"

$ws.Range("A21").Value = "wide:FlowModels"
$ws.Range("B21").Value = "Flow Models"
$ws.Range("D21").Value = ""
$ws.Range("F21").Value = ""

$ws.Range("A22").Value = "wide:AerolasticModels"
$ws.Range("B22").Value = "Aerolastic Models"
$ws.Range("D22").Value = ""
$ws.Range("F22").Value = ""

$ws.Range("A23").Value = "wide:ElectricalModels"
$ws.Range("B23").Value = "Electrical Models"
$ws.Range("D23").Value = ""
$ws.Range("F23").Value = ""

$ws.Range("A24").Value = "wide:FinancialModels"
$ws.Range("B24").Value = "Financial Models"
$ws.Range("D24").Value = ""
$ws.Range("F24").Value = ""

$ws.Range("A25").Value = "wide:GCM"
$ws.Range("B25").Value = "GCM"
$ws.Range("D25").Value = ""
$ws.Range("F25").Value = "wide:MeteorologicalModels"

$ws.Range("A26").Value = "wide:Mesoscale"
$ws.Range("B26").Value = "Mesoscale"
$ws.Range("D26").Value = ""
$ws.Range("F26").Value = "wide:MeteorologicalModels"

$ws.Range("A27").Value = "wide:Hindcast"
$ws.Range("B27").Value = "Hindcast"
$ws.Range("D27").Value = ""
$ws.Range("F27").Value = "wide:MeteorologicalModels"

$ws.Range("A28").Value = "wide:MorisonEquation"
$ws.Range("B28").Value = "Morison Equation"
$ws.Range("D28").Value = ""
$ws.Range("F28").Value = "wide:HydrodynamicModels"

$ws.Range("A29").Value = "wide:Radiation-Diffraction"
$ws.Range("B29").Value = "Radiation-Diffraction"
$ws.Range("D29").Value = ""
$ws.Range("F29").Value = "wide:HydrodynamicModels"

$ws.Range("A30").Value = "wide:Linear"
$ws.Range("B30").Value = "Linear"
$ws.Range("D30").Value = ""
$ws.Range("F30").Value = "wide:WaveModels"

$ws.Range("A31").Value = "wide:SecondOrder"
$ws.Range("B31").Value = "Second Order"
$ws.Range("D31").Value = ""
$ws.Range("F31").Value = "wide:WaveModels"

$ws.Range("A32").Value = "wide:FullyNonlinear"
$ws.Range("B32").Value = "Fully Nonlinear"
$ws.Range("D32").Value = ""
$ws.Range("F32").Value = "wide:WaveModels"

$ws.Range("A33").Value = "wide:Linearized"
$ws.Range("B33").Value = "Linearized"
$ws.Range("D33").Value = ""
$ws.Range("F33").Value = "wide:FlowModels"

$ws.Range("A34").Value = "wide:RANS"
$ws.Range("B34").Value = "RANS"
$ws.Range("D34").Value = ""
$ws.Range("F34").Value = "wide:FlowModels"

$ws.Range("A35").Value = "wide:LES"
$ws.Range("B35").Value = "LES"
$ws.Range("D35").Value = ""
$ws.Range("F35").Value = "wide:FlowModels"

$ws.Range("A36").Value = "wide:DNS"
$ws.Range("B36").Value = "DNS"
$ws.Range("D36").Value = ""
$ws.Range("F36").Value = "wide:FlowModels"

$ws.Range("A37").Value = "wide:VortexMethod"
$ws.Range("B37").Value = "Vortex Method"
$ws.Range("D37").Value = ""
$ws.Range("F37").Value = "wide:FlowModels"

$ws.Range("A38").Value = "wide:Analytical"
$ws.Range("B38").Value = "Analytical"
$ws.Range("D38").Value = ""
$ws.Range("F38").Value = "wide:FlowModels"

$ws.Range("A39").Value = "wide:Physical"
$ws.Range("B39").Value = "Physical"
$ws.Range("D39").Value = ""
$ws.Range("F39").Value = "wide:FlowModels"

$ws.Range("A40").Value = "wide:FEM"
$ws.Range("B40").Value = "FEM"
$ws.Range("D40").Value = ""
$ws.Range("F40").Value = "wide:AerolasticModels"

$ws.Range("A41").Value = "wide:PowerFlow"
$ws.Range("B41").Value = "Power Flow"
$ws.Range("D41").Value = ""
$ws.Range("F41").Value = "wide:ElectricalModels"

$ws.Range("A42").Value = "wide:OPF"
$ws.Range("B42").Value = "OPF"
$ws.Range("D42").Value = ""
$ws.Range("F42").Value = "wide:ElectricalModels"

$ws.Range("A43").Value = "wide:Small-SignalModels"
$ws.Range("B43").Value = "Small-Signal Models"
$ws.Range("D43").Value = ""
$ws.Range("F43").Value = "wide:ElectricalModels"

$ws.Range("A44").Value = "wide:DynamicModels"
$ws.Range("B44").Value = "Dynamic Models"
$ws.Range("D44").Value = ""
$ws.Range("F44").Value = "wide:ElectricalModels"

$ws.Range("A45").Value = "wide:ShortCircuitModels"
$ws.Range("B45").Value = "Short Circuit Models"
$ws.Range("D45").Value = ""
$ws.Range("F45").Value = "wide:ElectricalModels"

$ws.Range("A46").Value = "wide:StateEstimation"
$ws.Range("B46").Value = "State Estimation"
$ws.Range("D46").Value = ""
$ws.Range("F46").Value = "wide:ElectricalModels"

$ws.Range("A47").Value = "wide:PowerProtectionAnalysisModels"
$ws.Range("B47").Value = "Power Protection Analysis Models"
$ws.Range("D47").Value = ""
$ws.Range("F47").Value = "wide:ElectricalModels"

$ws.Range("A48").Value = "wide:ContingencyAnalysisModels"
$ws.Range("B48").Value = "Contingency Analysis Models"
$ws.Range("D48").Value = ""
$ws.Range("F48").Value = "wide:ElectricalModels"

$ws.Range("A49").Value = "wide:HarmonicModels"
$ws.Range("B49").Value = "Harmonic Models"
$ws.Range("D49").Value = ""
$ws.Range("F49").Value = "wide:ElectricalModels"

$ws.Range("A50").Value = "wide:BoS"
$ws.Range("B50").Value = "BoS"
$ws.Range("D50").Value = ""
$ws.Range("F50").Value = "wide:FinancialModels"

$ws.Range("A51").Value = "wide:NPV"
$ws.Range("B51").Value = "NPV"
$ws.Range("D51").Value = ""
$ws.Range("F51").Value = "wide:FinancialModels"

$ws.Range("A52").Value = "wide:LCOE"
$ws.Range("B52").Value = "LCOE"
$ws.Range("D52").Value = ""
$ws.Range("F52").Value = "wide:FinancialModels"

$ws.Range("A53").Value = "wide:IRR"
$ws.Range("B53").Value = "IRR"
$ws.Range("D53").Value = ""
$ws.Range("F53").Value = "wide:FinancialModels"

$ws.Range("A54").Value = "wide:"
$ws.Range("B54").Value = ""
$ws.Range("D54").Value = ""
$ws.Range("F54").Value = ""

